$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H (roboticRNAPrep) holds boolean FALSE values (displayed as "FALSE" via a custom
# number format). Convert them to literal text string "False" for rows 2-27.
for ($r = 2; $r -le 27; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $cell.NumberFormat = "@"
    $cell.Value = "False"
}

# Update the active selection to match the authored state.
$ws.Range("J15").Select()
